# Applies the two changes from the commit:
#  1. Slide 1, "TextBox 4": merge the two runs
#       "类中存在与基类中函数相同的函数" + "名，并满足如下两种情况："
#     into a single run of text (same run formatting on both, so this is a
#     pure text/run-structure simplification, no visible content change).
#  2. Slide 1, "TextBox 5": reposition/resize the shape from
#       off  (4821715, 3120057) ext (3312368, 3139321)
#     to
#       off  (4139952, 3120057) ext (4320480, 2862322)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Change 1: merge the two adjacent runs into one ------------------------
$body = $s.Shapes.Item(2)
$tr = $body.TextFrame.TextRange
# The two runs being merged start at character 80 of the shape's text
# (paragraph "子类中存在与基类中函数相同的函数名，并满足如下两种情况：...")
# and together span 27 characters ("类中存在与基类中函数相同的函数" = 15
# chars + "名，并满足如下两种情况：" = 12 chars).
$run = $tr.Characters(80, 27)
$run.Text = "类中存在与基类中函数相同的函数名，并满足如下两种情况："

# Re-assigning text inside an auto-fit ("Resize shape to fit text") textbox
# recomputes the shape's fitted height as a side effect, even though the
# overall text content/length here is unchanged. Put the box's size back to
# its original value so only the intended run-merge is reflected.
$body.Width = 623.6913452148438
$body.Height = 312.62347412109375

# --- Change 2: move/resize "TextBox 5" --------------------------------------
$box = $s.Shapes.Item(3)
# Values below are the EMU targets (4139952, 3120057, 4320480, 2862322)
# expressed in points (EMU / 12700), chosen so that the host's float32
# Left/Top/Width/Height storage reproduces the exact target EMU values.
$box.Left = 325.9804992675781
$box.Top = 245.6737823486328
$box.Width = 340.1952819824219
$box.Height = 225.37969970703125
